# Regenerate the handback-status report with a new pair of file GUIDs and
# fresh "Generate Report" timestamps.
#
#   1aa29009-39e0-4b33-a645-3f348e20e891  ->  fdc9c1bb-7730-40ed-81eb-6bcf496919e7
#   603718cb-1111-4a69-ba0a-989b0d347a7d  ->  ffff85af2807-2b9e-4e92-ab25-8fd6df1e79ac
#   1d419a78037f0c5f01dfa176c821250c8473c753 -> 76e99ae8dbfe0cc7975c923e647324dfd4707f28
#   57328d7613f7bf05c785a2af73361c52d54d9c34 -> 76e99ae8dbfe0cc7975c923e647324dfd4707f28
#
# plus the timestamps that moved forward by ~1m12s on the new CI run.

$wb = $excel.ActiveWorkbook

# old -> new: 1aa29009-39e0-4b33-a645-3f348e20e891 -> fdc9c1bb-7730-40ed-81eb-6bcf496919e7
$guid1New = "fdc9c1bb-7730-40ed-81eb-6bcf496919e7"
# old -> new: 603718cb-1111-4a69-ba0a-989b0d347a7d -> ffff85af2807-2b9e-4e92-ab25-8fd6df1e79ac
$guid2New = "ffff85af2807-2b9e-4e92-ab25-8fd6df1e79ac"

# old -> new: 1d419a78037f0c5f01dfa176c821250c8473c753 / 57328d7613f7bf05c785a2af73361c52d54d9c34
#             -> 76e99ae8dbfe0cc7975c923e647324dfd4707f28 (both collapse to the same new hash)
$hash1New = "76e99ae8dbfe0cc7975c923e647324dfd4707f28"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "$guid1New.md"
$ws.Range("B2").Value = "e2e\$guid1New.md"
$ws.Range("G2").Value = "2016-08-26 17:03:35"

$ws.Range("A3").Value = "$guid2New.md"
$ws.Range("B3").Value = "e2e\$guid2New.md"
$ws.Range("G3").Value = "2016-08-26 17:03:35"

foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$B$2') {
        $h.TextToDisplay = "e2e\$guid1New.md"
    }
    if ($addr -eq '$B$3') {
        $h.TextToDisplay = "e2e\$guid2New.md"
    }
}

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "$guid1New.md"
$ws.Range("G2").Value = "$guid1New.$hash1New.zh-cn.xlf"
$ws.Range("H2").Value = "2016-08-26 17:03:30"
$ws.Range("I2").Value = "$guid1New.md"
$ws.Range("J2").Value = "$guid1New.$hash1New.zh-cn.xlf"
$ws.Range("K2").Value = "2016-08-26 17:03:47"

$ws.Range("A3").Value = "$guid2New.md"
$ws.Range("G3").Value = "$guid1New.$hash1New.zh-cn.xlf"
$ws.Range("H3").Value = "2016-08-26 17:03:30"
$ws.Range("I3").Value = "$guid2New.md"
$ws.Range("J3").Value = "$guid1New.$hash1New.zh-cn.xlf"
$ws.Range("K3").Value = "2016-08-26 17:03:47"

foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = "$guid1New.md"
    }
    if ($addr -eq '$I$2') {
        $h.TextToDisplay = "$guid1New.md"
    }
    if ($addr -eq '$A$3') {
        $h.TextToDisplay = "$guid2New.md"
    }
    if ($addr -eq '$I$3') {
        $h.TextToDisplay = "$guid2New.md"
    }
}

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "$guid1New.md"
$ws.Range("G2").Value = "$guid1New.$hash1New.de-de.xlf"
$ws.Range("H2").Value = "2016-08-26 17:03:35"
$ws.Range("I2").Value = "$guid1New.md"
$ws.Range("J2").Value = "$guid1New.$hash1New.de-de.xlf"
$ws.Range("K2").Value = "2016-08-26 17:03:55"

$ws.Range("A3").Value = "$guid2New.md"
$ws.Range("G3").Value = "$guid1New.$hash1New.de-de.xlf"
$ws.Range("H3").Value = "2016-08-26 17:03:35"
$ws.Range("I3").Value = "$guid2New.md"
$ws.Range("J3").Value = "$guid1New.$hash1New.de-de.xlf"
$ws.Range("K3").Value = "2016-08-26 17:03:55"

foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = "$guid1New.md"
    }
    if ($addr -eq '$I$2') {
        $h.TextToDisplay = "$guid1New.md"
    }
    if ($addr -eq '$A$3') {
        $h.TextToDisplay = "$guid2New.md"
    }
    if ($addr -eq '$I$3') {
        $h.TextToDisplay = "$guid2New.md"
    }
}
